# Apply the "Finished Order-Settings-Design, Added Trendstrength Weighting
# and Updated ZeitNoctua_Zeitaufzeichnung" edit to the BacktestingSoftware sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BacktestingSoftware")

# --- Step 1: "finish" the old Order-Settings-Screen entry (row 23, col F) -------
# In the original file F23 held "Order-Settings-Screen 20%". The author
# repurposed that same work-log line to describe the newly finished
# trend-strength weighting work, then logged the (new) Order-Settings-Screen
# 40% status as a fresh line. Re-creating that edit sequence exactly keeps the
# shared-string table layout identical to the authored workbook.
$ws.Range("F23").Value = "Verschiedene Trendstärken in die Performancemessung integriert"
$movedText = $ws.Range("F23").Value2

# --- Step 2: add the two new work-log rows ---------------------------------
$ws.Range("A24").Value = "Pawlowsky"
$ws.Range("B24").Value = "Performancemessung integrieren"
$ws.Range("C24").Value = $movedText
$ws.Range("D24").Value = 41291
$ws.Range("E24").Value = 6

$ws.Range("A25").Value = "Nagy"
$ws.Range("B25").Value = "Performancemessung integrieren"
$ws.Range("C25").Value = "Trendstärken-Berechnungsmodell erstellt"
$ws.Range("D25").Value = 41291
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = "Berechnungsmodell erstellt"

$ws.Range("F24").Value = "Trendstärken integriert 80%"

# --- Step 3: set the final (new) status text for the Order-Settings-Screen row -
$ws.Range("F23").Value = "Order-Settings-Screen 40%"

# --- Step 4: match formatting of the appended rows to the existing table ----
$ws.Range("A23:F23").Copy() | Out-Null
$ws.Range("A24:F24").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A23:F23").Copy() | Out-Null
$ws.Range("A25:F25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Step 5: restore the sheet selection the author left the sheet in -------
$ws.Activate()
$ws.Range("B12:F12").Select() | Out-Null

$wb.Save()
